# Refresh the crypto price/volume snapshot (Price = column D, Volume(1h) = column E)
# to match the latest scrape, per the "Updated cryptos list ... with GitHub Actions" commit.
#
# Row 40/41 additionally swap places in the ranking (MXToken now ranks above ARBITRUM),
# so every column (Coin/Link/Price/Volume) is rewritten for those two rows.
#
# Some "Price" values are plain decimals (e.g. 0.986) that Excel would otherwise
# auto-coerce to a number on assignment; prefixing with a leading apostrophe forces
# them to stay text, matching the original inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Leading apostrophe = "treat as text" quote-prefix, exactly like typing it in the
    # Excel UI; keeps numeric-looking strings (e.g. "0.986") from becoming real numbers.
    $ws.Range($range).Value = "'" + $text
}

$ws.Range("D2").Value = '27.715.36'
$ws.Range("E2").Value = '  +2.23%  '
$ws.Range("D3").Value = '1.573.42'
$ws.Range("E3").Value = '  +0.26%  '
Set-TextValue "D4" '0.986'
$ws.Range("E4").Value = '  -2.39%  '
Set-TextValue "D5" '211.19'
$ws.Range("E5").Value = '  +0.45%  '
Set-TextValue "D6" '0.494'
$ws.Range("E6").Value = '  +0.53%  '
Set-TextValue "D7" '0.987'
$ws.Range("E7").Value = '  -1.91%  '
Set-TextValue "D8" '23.32'
$ws.Range("E8").Value = '  +5.66%  '
$ws.Range("E9").Value = '  +0.90%  '
Set-TextValue "D10" '0.0598'
$ws.Range("E10").Value = '  +0.10%  '
Set-TextValue "D11" '0.0875'
$ws.Range("E11").Value = '  +1.49%  '
$ws.Range("D12").Value = '1.793.02'
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("D13").Value = '1.589.91'
$ws.Range("E13").Value = '  +1.34%  '
Set-TextValue "D14" '3.77'
$ws.Range("E14").Value = '  -0.28%  '
Set-TextValue "D15" '0.522'
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").Value = '27.604.08'
$ws.Range("E16").Value = '  +1.95%  '
Set-TextValue "D17" '63.22'
$ws.Range("E17").Value = '  +1.79%  '
Set-TextValue "D18" '230.49'
$ws.Range("E18").Value = '  +6.51%  '
Set-TextValue "D19" '7.52'
$ws.Range("E19").Value = '  +1.55%  '
$ws.Range("D20").Value = '0.0₃0706'
$ws.Range("E20").Value = '  +0.14%  '
Set-TextValue "D21" '0.988'
$ws.Range("E21").Value = '  -1.97%  '
Set-TextValue "D22" '4.13'
$ws.Range("E22").Value = '  -0.58%  '
Set-TextValue "D23" '9.44'
$ws.Range("E23").Value = '  +2.55%  '
Set-TextValue "D24" '1.97'
$ws.Range("E24").Value = '  +1.21%  '
Set-TextValue "D25" '149.84'
$ws.Range("E25").Value = '  -2.57%  '
Set-TextValue "D26" '15.33'
$ws.Range("E26").Value = '  +1.30%  '
Set-TextValue "D27" '6.60'
$ws.Range("E27").Value = '  -0.17%  '
$ws.Range("E28").Value = '  +1.59%  '
$ws.Range("E29").Value = '  -1.49%  '
$ws.Range("E30").Value = '  +0.05%  '
Set-TextValue "D31" '0.0474'
$ws.Range("E31").Value = '  +0.07%  '
Set-TextValue "D32" '3.25'
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("D33").Value = '1.454.74'
$ws.Range("E33").Value = '  +1.24%  '
Set-TextValue "D34" '3.13'
$ws.Range("E34").Value = '  -1.73%  '
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("E36").Value = '  -6.05%  '
$ws.Range("E37").Value = '  -0.98%  '
$ws.Range("E38").Value = '  +0.85%  '
Set-TextValue "D39" '0.543'
$ws.Range("E39").Value = '  +2.00%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D40" '2.46'
$ws.Range("E40").Value = '  +3.55%  '
$ws.Range("B41").Value = 'ARBITRUM'
$ws.Range("C41").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D41" '0.814'
$ws.Range("E41").Value = '  +0.18%  '
$ws.Range("E42").Value = '  -3.17%  '
$ws.Range("E43").Value = '  -2.00%  '
Set-TextValue "D44" '1.86'
$ws.Range("E44").Value = '  +7.09%  '
$ws.Range("E45").Value = '  -2.81%  '
Set-TextValue "D46" '64.26'
$ws.Range("E46").Value = '  -0.79%  '
$ws.Range("D47").Value = '1.705.80'
Set-TextValue "D48" '86.98'
$ws.Range("E48").Value = '  +1.53%  '
Set-TextValue "D49" '0.0524'
$ws.Range("E49").Value = '  +1.22%  '
$ws.Range("D50").Value = '0.0₇0995'
$ws.Range("E50").Value = '  -2.43%  '
Set-TextValue "D51" '40.32'
$ws.Range("E51").Value = '  +18.68%  '
